$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 37: trade finished -> status "DONE", finalized date filled in, fee recorded ---
$ws.Cells.Item(37, 8).Value = "DONE"
$ws.Cells.Item(37, 9).Value = 42858.746388888889
$ws.Cells.Item(37, 10).Value = "0.00147494 USDT (0.15%)"

# --- Row 38: trade cancelled -> status "CANCEL", finalized date filled in ---
$ws.Cells.Item(38, 8).Value = "CANCEL"
$ws.Cells.Item(38, 9).Value = 42859.441886574074

# --- Row 39: brand new trade entry (ETC/USDT buy, still IN PROGRESS) ---
# Clone formatting from the row above first (cell by cell) so number formats /
# wrap settings match the rest of the table, then fill in the values.
$ws.Cells.Item(38, 1).Copy() | Out-Null
$ws.Cells.Item(39, 1).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(39, 1).Value = 42859.441886574074

$ws.Cells.Item(39, 2).Value = "            Buy"
$ws.Cells.Item(39, 3).Value = "        ETC"

$ws.Cells.Item(39, 4).Value = "'                7.64000000`n"
$ws.Cells.Item(38, 4).Copy() | Out-Null
$ws.Cells.Item(39, 4).PasteSpecial(-4122) | Out-Null

$ws.Cells.Item(39, 5).Value = "          7.1USDT"
$ws.Cells.Item(39, 6).Value = "        3.4ETC"
$ws.Cells.Item(39, 7).Value = " ETC/USDT0000001"
$ws.Cells.Item(39, 8).Value = "IN PROGRESS"

$ws.Cells.Item(38, 9).Copy() | Out-Null
$ws.Cells.Item(39, 9).PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = $false

# --- Update the active selection to reflect where editing left off ---
$ws.Range("B43").Select()
